$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 64, shifting rows 64:70 down to 65:71
$ws.Rows.Item(64).Insert()

# Populate the new row 64 with data (new record)
$ws.Cells.Item(64, 1).Value = 7
$ws.Cells.Item(64, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(64, 3).Value = "Ñuble"
$ws.Cells.Item(64, 4).Value = 44610
$ws.Cells.Item(64, 4).NumberFormat = $ws.Cells.Item(65, 4).NumberFormat
$ws.Cells.Item(64, 5).Value = 16
$ws.Cells.Item(64, 6).Value = 100112031
$ws.Cells.Item(64, 7).Value = "Poroto verde"
$ws.Cells.Item(64, 8).Value = "Sin especificar"
$ws.Cells.Item(64, 9).Value = "Primera"
$ws.Cells.Item(64, 10).Value = 100
$ws.Cells.Item(64, 11).Value = 29000
$ws.Cells.Item(64, 12).Value = 30000
$ws.Cells.Item(64, 13).Value = 29500
$ws.Cells.Item(64, 14).Value = "$/saco 25 kilos"
$ws.Cells.Item(64, 15).Value = "Región del Maule"
$ws.Cells.Item(64, 16).Value = 1180
$ws.Cells.Item(64, 17).Value = 25
$ws.Cells.Item(64, 18).Value = "Hortaliza"
